$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.128.87"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.32"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.36"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4754"
$ws.Range("E7").Value = "  -2.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2815"
$ws.Range("E8").Value = "  -3.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06467"
$ws.Range("E9").Value = "  -3.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.864.88"
$ws.Range("E10").Value = "  -1.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07279"
$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.27"
$ws.Range("E12").Value = "  -4.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.119"
$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.06"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6439"
$ws.Range("E15").Value = "  -3.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.080.58"
$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.22"
$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007612"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.096.33"
$ws.Range("E20").Value = "  -2.55%  "

$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.239"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "216.65"
$ws.Range("E23").Value = "  +13.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.093"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.217"
$ws.Range("E25").Value = "  -2.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.47"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.913"
$ws.Range("E28").Value = "  -0.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.428"
$ws.Range("E29").Value = "  -2.37%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09152"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.222"
$ws.Range("E31").Value = "  -2.90%  "

$ws.Range("E32").Value = "  -3.00%  "

$ws.Range("E33").Value = "  -3.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7403"
$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  +3.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.685"
$ws.Range("E36").Value = "  -1.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01815"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.606"
$ws.Range("E38").Value = "  -2.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8992"
$ws.Range("E39").Value = "  -1.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.044"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.919"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.00"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.0000"
$ws.Range("E43").Value = "  +0.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4228"
$ws.Range("E44").Value = "  -3.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.396"
$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("E46").Value = "  -5.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.549"
$ws.Range("E47").Value = "  +9.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.70"
$ws.Range("E48").Value = "  -7.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.776"
$ws.Range("E49").Value = "  -2.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.11"
$ws.Range("E50").Value = "  -2.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05681"
$ws.Range("E51").Value = "  -2.56%  "
